# Auto-generated Excel COM-interop script
# Applies the 'scheduled runner' price/profit refresh to the Leve-tracking sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). Columns touched per row:
#   H  currentAveragePrice
#   I  currentAveragePriceNQ
#   J  currentAveragePriceHQ
#   K  LevePriceNQ
#   L  LevePriceHQ
#   M  LeveProfitNQ
#   N  LeveProfitHQ

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 49056.43
$ws.Range("I33").Value = 67013.53
$ws.Range("K33").Value = 67013.53
$ws.Range("M33").Value = -66784.53

$ws.Range("H62").Value = 2162.6667
$ws.Range("I62").Value = 2162.6667
$ws.Range("K62").Value = 2162.6667
$ws.Range("M62").Value = -1538.6667

$ws.Range("H65").Value = 2162.6667
$ws.Range("I65").Value = 2162.6667
$ws.Range("K65").Value = 10813.3335
$ws.Range("M65").Value = -7693.333500000001

$ws.Range("H106").Value = 4959.8
$ws.Range("I106").Value = 5449.75
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 5449.75
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -4818.75
$ws.Range("N106").Value = -4262

$ws.Range("H132").Value = 5686780
$ws.Range("I132").Value = 6762484
$ws.Range("K132").Value = 20287452
$ws.Range("M132").Value = -20284922

$ws.Range("H141").Value = 1368.8226
$ws.Range("I141").Value = 1286.1167
$ws.Range("K141").Value = 3858.3501
$ws.Range("M141").Value = 1321.6499


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21176.396
$ws.Range("I32").Value = 4676.4873
$ws.Range("J32").Value = 141175.73
$ws.Range("K32").Value = 4676.4873
$ws.Range("L32").Value = 141175.73
$ws.Range("M32").Value = -4389.4873
$ws.Range("N32").Value = -141749.73

$ws.Range("H61").Value = 1009.38464
$ws.Range("I61").Value = 849.3111
$ws.Range("J61").Value = 2038.4286
$ws.Range("K61").Value = 849.3111
$ws.Range("L61").Value = 2038.4286
$ws.Range("M61").Value = -637.3111
$ws.Range("N61").Value = -2462.4286

$ws.Range("H74").Value = 459.3158
$ws.Range("I74").Value = 428.79413
$ws.Range("J74").Value = 718.75
$ws.Range("K74").Value = 428.79413
$ws.Range("L74").Value = 718.75
$ws.Range("M74").Value = 445.20587
$ws.Range("N74").Value = -2466.75

$ws.Range("H77").Value = 459.3158
$ws.Range("I77").Value = 428.79413
$ws.Range("J77").Value = 718.75
$ws.Range("K77").Value = 2143.97065
$ws.Range("L77").Value = 3593.75
$ws.Range("M77").Value = 2224.02935
$ws.Range("N77").Value = -12329.75

$ws.Range("H110").Value = 71579130
$ws.Range("I110").Value = 77085140
$ws.Range("K110").Value = 77085140
$ws.Range("M110").Value = -77083095

$ws.Range("H132").Value = 3765.75
$ws.Range("I132").Value = 3721.525
$ws.Range("J132").Value = 3986.875
$ws.Range("K132").Value = 11164.575
$ws.Range("L132").Value = 11960.625
$ws.Range("M132").Value = -8634.575000000001
$ws.Range("N132").Value = -17020.625

$ws.Range("H136").Value = 1009.38464
$ws.Range("I136").Value = 849.3111
$ws.Range("J136").Value = 2038.4286
$ws.Range("K136").Value = 2547.9333
$ws.Range("L136").Value = 6115.2858
$ws.Range("M136").Value = 2.066699999999855
$ws.Range("N136").Value = -11215.2858


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H52").Value = 45819.25
$ws.Range("I52").Value = 50592.332
$ws.Range("J52").Value = 31500
$ws.Range("K52").Value = 50592.332
$ws.Range("L52").Value = 31500
$ws.Range("M52").Value = -50329.332
$ws.Range("N52").Value = -32026

$ws.Range("H120").Value = 38333.332
$ws.Range("J120").Value = 38333.332
$ws.Range("L120").Value = 38333.332
$ws.Range("N120").Value = -48009.332

$ws.Range("H121").Value = 45819.25
$ws.Range("I121").Value = 50592.332
$ws.Range("J121").Value = 31500
$ws.Range("K121").Value = 50592.332
$ws.Range("L121").Value = 31500
$ws.Range("M121").Value = -48845.332
$ws.Range("N121").Value = -34994


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1421.25
$ws.Range("I16").Value = 979.1667
$ws.Range("K16").Value = 979.1667
$ws.Range("M16").Value = -692.1667

$ws.Range("H31").Value = 26670.207
$ws.Range("I31").Value = 1513.4688
$ws.Range("J31").Value = 57632.348
$ws.Range("K31").Value = 1513.4688
$ws.Range("L31").Value = 57632.348
$ws.Range("M31").Value = -1218.4688
$ws.Range("N31").Value = -58222.348

$ws.Range("H34").Value = 26670.207
$ws.Range("I34").Value = 1513.4688
$ws.Range("J34").Value = 57632.348
$ws.Range("K34").Value = 1513.4688
$ws.Range("L34").Value = 57632.348
$ws.Range("M34").Value = -1311.4688
$ws.Range("N34").Value = -58036.348

$ws.Range("H58").Value = 1232.98
$ws.Range("I58").Value = 1034.3489
$ws.Range("J58").Value = 2453.1428
$ws.Range("K58").Value = 1034.3489
$ws.Range("L58").Value = 2453.1428
$ws.Range("M58").Value = -831.3489
$ws.Range("N58").Value = -2859.1428

$ws.Range("H94").Value = 1028.7059
$ws.Range("I94").Value = 800.8
$ws.Range("K94").Value = 800.8
$ws.Range("M94").Value = -349.8

$ws.Range("H113").Value = 1421.25
$ws.Range("I113").Value = 979.1667
$ws.Range("K113").Value = 979.1667
$ws.Range("M113").Value = 1190.8333

$ws.Range("H132").Value = 3041.12
$ws.Range("I132").Value = 2826.7222
$ws.Range("J132").Value = 3592.4285
$ws.Range("K132").Value = 8480.1666
$ws.Range("L132").Value = 10777.2855
$ws.Range("M132").Value = -5950.1666
$ws.Range("N132").Value = -15837.2855

$ws.Range("H136").Value = 1232.98
$ws.Range("I136").Value = 1034.3489
$ws.Range("J136").Value = 2453.1428
$ws.Range("K136").Value = 3103.0467
$ws.Range("L136").Value = 7359.428400000001
$ws.Range("M136").Value = -553.0466999999999
$ws.Range("N136").Value = -12459.4284

$ws.Range("H140").Value = 4927.25
$ws.Range("I140").Value = 4927.25
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 4927.25
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = 252.75
$ws.Range("N140").Value = $null


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 99.333336
$ws.Range("I6").Value = 84.85714
$ws.Range("K6").Value = 254.57142
$ws.Range("M6").Value = -141.57142

$ws.Range("H10").Value = 400
$ws.Range("I10").Value = 200
$ws.Range("J10").Value = 600
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 1800
$ws.Range("M10").Value = -461
$ws.Range("N10").Value = -2078

$ws.Range("H122").Value = 410
$ws.Range("I122").Value = 416
$ws.Range("J122").Value = 350
$ws.Range("K122").Value = 3744
$ws.Range("L122").Value = 3150
$ws.Range("M122").Value = -1294
$ws.Range("N122").Value = -8050

$ws.Range("H132").Value = 1432.5
$ws.Range("I132").Value = 794
$ws.Range("J132").Value = 2071
$ws.Range("K132").Value = 7146
$ws.Range("L132").Value = 18639
$ws.Range("M132").Value = -4616
$ws.Range("N132").Value = -23699


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 83421144
$ws.Range("I80").Value = 200207360
$ws.Range("J80").Value = 2425.5715
$ws.Range("K80").Value = 200207360
$ws.Range("L80").Value = 2425.5715
$ws.Range("M80").Value = -200206362
$ws.Range("N80").Value = -4421.5715

$ws.Range("H83").Value = 83421144
$ws.Range("I83").Value = 200207360
$ws.Range("J83").Value = 2425.5715
$ws.Range("K83").Value = 1001036800
$ws.Range("L83").Value = 12127.8575
$ws.Range("M83").Value = -1001031808
$ws.Range("N83").Value = -22111.8575

$ws.Range("H102").Value = 2228.1785
$ws.Range("I102").Value = 1714.6
$ws.Range("J102").Value = 2820.7693
$ws.Range("K102").Value = 1714.6
$ws.Range("L102").Value = 2820.7693
$ws.Range("M102").Value = -92.59999999999991
$ws.Range("N102").Value = -6064.7693


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2369.5
$ws.Range("I93").Value = 2200.4167
$ws.Range("J93").Value = 2876.75
$ws.Range("K93").Value = 2200.4167
$ws.Range("L93").Value = 2876.75
$ws.Range("M93").Value = -952.4167000000002
$ws.Range("N93").Value = -5372.75

$ws.Range("H122").Value = 1986.1111
$ws.Range("I122").Value = 1893.0714
$ws.Range("J122").Value = 2311.75
$ws.Range("K122").Value = 5679.2142
$ws.Range("L122").Value = 6935.25
$ws.Range("M122").Value = -3229.2142
$ws.Range("N122").Value = -11835.25


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 510.3231
$ws.Range("I136").Value = 331.92307
$ws.Range("J136").Value = 1223.9231
$ws.Range("K136").Value = 995.7692099999999
$ws.Range("L136").Value = 3671.7693
$ws.Range("M136").Value = 1554.23079
$ws.Range("N136").Value = -8771.7693


Write-Output "Applied scheduled price/profit refresh across all 8 sheets."